# The workbook originally compared three routing algorithms side by side
# (Q_ROUTING, DIJKSTRA, BELLMAN_FORD), each on its own sheet. After the
# refactor only the Bellman-Ford results are kept, re-run against the
# refactored classes, so the other two sheets are removed and the
# remaining sheet's metrics are refreshed with the new run's numbers.

$wb = $excel.ActiveWorkbook

$wb.Worksheets("Q_ROUTING").Delete() | Out-Null
$wb.Worksheets("DIJKSTRA").Delete() | Out-Null

$ws = $wb.Worksheets("BELLMAN_FORD")

# Refresh episode metrics from the latest simulation run.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 18
$ws.Range("D2").Value = 17
$ws.Range("E2").Value = $false
$ws.Range("F2").Value = 116
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "{}"
